$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.802.19"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "2.290.36"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "120.29"
$ws.Range("E5").Value = "  +6.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.59"
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.640"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.39"
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0943"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.22"
$ws.Range("E12").Value = "  +2.12%  "
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.63"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("E15").Value = "  +4.49%  "
$ws.Range("D16").Value = "2.635.73"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "2.288.68"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "43.743.63"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.10"
$ws.Range("E20").Value = "  +3.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.62"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.47"
$ws.Range("E22").Value = "  +1.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.46"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.70"
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("E25").Value = "  -4.66%  "
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.83"
$ws.Range("E27").Value = "  +2.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.75"
$ws.Range("E28").Value = "  +2.82%  "
$ws.Range("E29").Value = "  +0.35%  "
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.42"
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.78"
$ws.Range("E32").Value = "  +1.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0914"
$ws.Range("E33").Value = "  -1.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.82"
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("E35").Value = "  +2.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0382"
$ws.Range("E36").Value = "  +4.83%  "
$ws.Range("E37").Value = "  +2.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.99"
$ws.Range("E38").Value = "  +4.20%  "
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("E40").Value = "  +8.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.35"
$ws.Range("E41").Value = "  +3.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "75.15"
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("E43").Value = "  -1.63%  "
$ws.Range("E44").Value = "  -5.39%  "
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("E47").Value = "  +3.20%  "
$ws.Range("E48").Value = "  -2.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "73.39"
$ws.Range("E49").Value = "  +36.61%  "
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "102.18"
$ws.Range("E51").Value = "  +0.73%  "
